$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 163
    3  = 423
    4  = 12319
    6  = 138
    10 = 191
    11 = 446
    12 = 58
    16 = 366
    17 = 3280
    22 = 31
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
